$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data between row pairs (order of fixtures corrected) ---
# Row 58 <-> Row 59
$ws.Range("F58").Value = 'Betis'
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 'Cadiz CF'
$ws.Range("I58").Value = 1
$ws.Range("J58").Value = 1.87
$ws.Range("K58").Value = '05/09/2023 12:02'
$ws.Range("L58").Value = 1.81
$ws.Range("M58").Value = '24/09/2023 18:26'
$ws.Range("N58").Value = 3.56
$ws.Range("O58").Value = '05/09/2023 12:02'
$ws.Range("P58").Value = 3.66
$ws.Range("Q58").Value = '24/09/2023 18:26'
$ws.Range("R58").Value = 4.54
$ws.Range("S58").Value = '05/09/2023 12:02'
$ws.Range("T58").Value = 5.07
$ws.Range("U58").Value = '24/09/2023 18:26'
$ws.Range("V58").Value = 'https://www.betexplorer.com/football/spain/laliga/betis-cadiz/IicoJIZo/'
$ws.Range("F59").Value = 'Las Palmas'
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 'Granada CF'
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2.14
$ws.Range("K59").Value = '11/09/2023 13:19'
$ws.Range("L59").Value = 1.95
$ws.Range("M59").Value = '24/09/2023 18:19'
$ws.Range("N59").Value = 3.21
$ws.Range("O59").Value = '11/09/2023 13:19'
$ws.Range("P59").Value = 3.73
$ws.Range("Q59").Value = '24/09/2023 18:27'
$ws.Range("R59").Value = 3.61
$ws.Range("S59").Value = '11/09/2023 13:19'
$ws.Range("T59").Value = 4.1
$ws.Range("U59").Value = '24/09/2023 18:27'
$ws.Range("V59").Value = 'https://www.betexplorer.com/football/spain/laliga/las-palmas-granada-cf/tWsBDE3N/'

# Row 64 <-> Row 65
$ws.Range("F64").Value = 'Real Madrid'
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 'Las Palmas'
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1.2
$ws.Range("K64").Value = '23/09/2023 09:28'
$ws.Range("L64").Value = 1.18
$ws.Range("M64").Value = '27/09/2023 18:29'
$ws.Range("N64").Value = 6.76
$ws.Range("O64").Value = '23/09/2023 09:28'
$ws.Range("P64").Value = 8
$ws.Range("Q64").Value = '27/09/2023 18:29'
$ws.Range("R64").Value = 11.3
$ws.Range("S64").Value = '23/09/2023 09:28'
$ws.Range("T64").Value = 16.5
$ws.Range("U64").Value = '27/09/2023 18:29'
$ws.Range("V64").Value = 'https://www.betexplorer.com/football/spain/laliga/real-madrid-las-palmas/GQHmRXXM/'
$ws.Range("F65").Value = 'Ath Bilbao'
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 'Getafe'
$ws.Range("I65").Value = 2
$ws.Range("J65").Value = 1.71
$ws.Range("K65").Value = '17/09/2023 09:02'
$ws.Range("L65").Value = 1.53
$ws.Range("M65").Value = '27/09/2023 18:31'
$ws.Range("N65").Value = 3.42
$ws.Range("O65").Value = '17/09/2023 09:02'
$ws.Range("P65").Value = 4.06
$ws.Range("Q65").Value = '27/09/2023 18:49'
$ws.Range("R65").Value = 5.44
$ws.Range("S65").Value = '17/09/2023 09:02'
$ws.Range("T65").Value = 7.73
$ws.Range("U65").Value = '27/09/2023 18:49'
$ws.Range("V65").Value = 'https://www.betexplorer.com/football/spain/laliga/ath-bilbao-getafe/zgsFCYIT/'

# Row 88 <-> Row 89
$ws.Range("F88").Value = 'Celta Vigo'
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 'Getafe'
$ws.Range("I88").Value = 2
$ws.Range("J88").Value = 1.92
$ws.Range("K88").Value = '28/09/2023 15:02'
$ws.Range("L88").Value = 2.04
$ws.Range("M88").Value = '08/10/2023 18:29'
$ws.Range("N88").Value = 3.25
$ws.Range("O88").Value = '28/09/2023 15:02'
$ws.Range("P88").Value = 3.29
$ws.Range("Q88").Value = '08/10/2023 18:27'
$ws.Range("R88").Value = 4.85
$ws.Range("S88").Value = '28/09/2023 15:02'
$ws.Range("T88").Value = 4.38
$ws.Range("U88").Value = '08/10/2023 18:29'
$ws.Range("V88").Value = 'https://www.betexplorer.com/football/spain/laliga/celta-vigo-getafe/0ARtdhXd/'
$ws.Range("F89").Value = 'Alaves'
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 'Betis'
$ws.Range("I89").Value = 1
$ws.Range("J89").Value = 2.8
$ws.Range("K89").Value = '01/10/2023 20:24'
$ws.Range("L89").Value = 2.57
$ws.Range("M89").Value = '08/10/2023 18:28'
$ws.Range("N89").Value = 3.04
$ws.Range("O89").Value = '01/10/2023 20:24'
$ws.Range("P89").Value = 3.25
$ws.Range("Q89").Value = '08/10/2023 18:28'
$ws.Range("R89").Value = 2.72
$ws.Range("S89").Value = '01/10/2023 20:24'
$ws.Range("T89").Value = 3.06
$ws.Range("U89").Value = '08/10/2023 18:22'
$ws.Range("V89").Value = 'https://www.betexplorer.com/football/spain/laliga/alaves-betis/YNPlfW19/'

# --- Append new fixture as row 138 (copy formatting of row 137 first) ---
$ws.Range("A137:V137").Copy()
$ws.Range("A138:V138").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A138").Value = 137
$ws.Range("B138").Value = 'spain'
$ws.Range("C138").Value = 'laliga'
$ws.Range("D138").Value = '2023-2024'
$ws.Range("E138").Value = 45256.875
$ws.Range("F138").Value = 'Betis'
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 'Las Palmas'
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 1.75
$ws.Range("K138").Value = '05/11/2023 11:03'
$ws.Range("L138").Value = 1.72
$ws.Range("M138").Value = '26/11/2023 20:37'
$ws.Range("N138").Value = 3.71
$ws.Range("O138").Value = '05/11/2023 11:03'
$ws.Range("P138").Value = 3.74
$ws.Range("Q138").Value = '26/11/2023 20:37'
$ws.Range("R138").Value = 4.97
$ws.Range("S138").Value = '05/11/2023 11:03'
$ws.Range("T138").Value = 5.55
$ws.Range("U138").Value = '26/11/2023 20:59'
$ws.Range("V138").Value = 'https://www.betexplorer.com/football/spain/laliga/betis-las-palmas/S0XoIWko/'
